$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.108.37'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '3.561.78'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '605.81'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').Value = '144.62'
$ws.Range('E6').Value = '  +0.58%  '
$ws.Range('D7').Value = '3.559.30'
$ws.Range('E7').Value = '  +1.87%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  +3.57%  '
$ws.Range('E10').Value = '  +0.69%  '
$ws.Range('D11').Value = '7.92'
$ws.Range('E11').Value = '  -2.54%  '
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('D13').Value = '4.163.89'
$ws.Range('E13').Value = '  +2.63%  '
$ws.Range('E14').Value = '  +2.04%  '
$ws.Range('D15').Value = '30.00'
$ws.Range('E15').Value = '  -0.56%  '
$ws.Range('D16').Value = '3.557.68'
$ws.Range('E16').Value = '  +3.20%  '
$ws.Range('D17').Value = '66.212.72'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('E18').Value = '  -0.59%  '
$ws.Range('D19').Value = '11.31'
$ws.Range('E19').Value = '  +7.80%  '
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('D21').Value = '14.84'
$ws.Range('E21').Value = '  +0.53%  '
$ws.Range('D22').Value = '429.10'
$ws.Range('E22').Value = '  +2.32%  '
$ws.Range('E23').Value = '  +4.44%  '
$ws.Range('D24').Value = '79.12'
$ws.Range('E24').Value = '  +2.08%  '
$ws.Range('D25').Value = '3.702.35'
$ws.Range('E25').Value = '  +2.89%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +1.56%  '
$ws.Range('E28').Value = '  +2.01%  '
$ws.Range('D29').Value = '7.95'
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('D30').Value = '9.08'
$ws.Range('E30').Value = '  -3.44%  '
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('D32').Value = '25.57'
$ws.Range('E32').Value = '  +2.20%  '
$ws.Range('E33').Value = '  -1.03%  '
$ws.Range('D34').Value = '3.556.02'
$ws.Range('E34').Value = '  +2.60%  '
$ws.Range('D35').Value = '0.153'
$ws.Range('E35').Value = '  -6.08%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('E37').Value = '  +1.78%  '
$ws.Range('E38').Value = '  +3.00%  '
$ws.Range('D39').Value = '5.61'
$ws.Range('E39').Value = '  +0.46%  '
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').Value = '176.00'
$ws.Range('E41').Value = '  +3.60%  '
$ws.Range('D42').Value = '0.0849'
$ws.Range('E42').Value = '  -2.41%  '
$ws.Range('D43').Value = '5.21'
$ws.Range('E43').Value = '  +1.80%  '
$ws.Range('D44').Value = '0.893'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = '1.95'
$ws.Range('D46').Value = '46.06'
$ws.Range('E46').Value = '  +0.92%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '25.76'
$ws.Range('E47').Value = '  -2.53%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value = '1.20'
$ws.Range('E48').Value = '  -1.25%  '
$ws.Range('D49').Value = '23.49'
$ws.Range('E49').Value = '  +8.49%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '7.12'
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').Value = '2.34'
$ws.Range('E51').Value = '  -0.19%  '
